# Brands.xlsx update:
#  - Sheet1: drop the two trailing blank rows (3 & 4) so the sheet's used
#    range shrinks back down to A1:B2.
#  - Sheet2: no explicit change needed here - it stops being the active/
#    selected tab once a different sheet becomes active below.
#  - Add a brand-new "Sheet3" after the existing sheets, put the value
#    "6787-897" in A1, and make it the active sheet (so it becomes the
#    selected tab, matching the new activeTab index in the workbook view).

$wb = $excel.ActiveWorkbook

# --- Sheet1: remove empty rows 3 and 4 -------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A3:A4").EntireRow.Delete() | Out-Null

# --- Add Sheet3 after the last existing sheet ------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"
$ws3.Range("A1").Value = "6787-897"

# Make the new sheet the active/selected tab.
$ws3.Activate()
